# Insert a new data row at row 92 (pushing existing rows 92-141 down to 93-142)
# and populate it with a new Mango price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(92).Insert()

$ws.Range("A92").Value = 7
$ws.Range("B92").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C92").Value = "Ñuble"
$ws.Range("D92").Value = 45086
$ws.Range("E92").Value = 16
$ws.Range("F92").Value = "Fruta"
$ws.Range("G92").Value = 100108
$ws.Range("H92").Value = "Tropicales y subtropicales"
$ws.Range("I92").Value = 100108002
$ws.Range("J92").Value = "Mango"
$ws.Range("K92").Value = "Sin especificar"
$ws.Range("L92").Value = "Primera"
$ws.Range("M92").Value = 60
$ws.Range("N92").Value = 9000
$ws.Range("O92").Value = 9000
$ws.Range("P92").Value = 9000
$ws.Range("Q92").Value = "$/bandeja 4 kilos"
$ws.Range("R92").Value = "Perú"
$ws.Range("S92").Value = 2250
$ws.Range("T92").Value = 4

# Match the date-format style used by the rest of column D (s="2")
$ws.Range("D92").NumberFormat = $ws.Range("D93").NumberFormat
